$d = $word.ActiveDocument
$ell = [char]8230

function Merge-WholeParagraphText($matchSnippet, $newText) {
    # Finds the paragraph whose text contains $matchSnippet and replaces the
    # entire paragraph's visible text (minus the paragraph mark) with $newText,
    # collapsing all runs/proofErr marks into a single run.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($matchSnippet)) {
            $rng = $p.Range
            $found = $rng.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
            return $found
        }
    }
    return $false
}

function Split-RunAt($paraMatchSnippet, $offsetIntoParagraph) {
    # Forces a run boundary inside a paragraph at a given character offset
    # (0-based, relative to the paragraph's Range.Start) without altering
    # formatting, by adding then immediately deleting a bookmark there.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($paraMatchSnippet)) {
            $full = $p.Range
            $boundary = $full.Start + $offsetIntoParagraph
            $pt = $d.Range($boundary, $boundary)
            $d.Bookmarks.Add("iron_tmp_split", $pt)
            $d.Bookmarks.Item("iron_tmp_split").Delete()
            return $true
        }
    }
    return $false
}

# ---------------------------------------------------------------------------
# 1) "Alberto " / "Perales" / " " -> "Alberto Perales " (single run)
# ---------------------------------------------------------------------------
Merge-WholeParagraphText "Perales" "Alberto Perales " | Out-Null

# ---------------------------------------------------------------------------
# 2) "b) " / "sub" / " goals ... seed " / "withing" / " pairs in boat " ->
#    single run "b) sub goals are to find right pairing of animals and seed withing pairs in boat "
# ---------------------------------------------------------------------------
Merge-WholeParagraphText "withing" "b) sub goals are to find right pairing of animals and seed withing pairs in boat " | Out-Null

# ---------------------------------------------------------------------------
# 3) "a) " stays its own run; "find" / " pairing to travel ....Cat..." ->
#    single run "find pairing to travel ….Cat and man , seed and man, parrot and man "
#    (the "a) " run must remain separate)
# ---------------------------------------------------------------------------
$mergedFind = "find pairing to travel " + $ell + ".Cat and man , seed and man, parrot and man "
Merge-WholeParagraphText "pairing to travel" ("a) " + $mergedFind) | Out-Null
Split-RunAt "pairing to travel" 3 | Out-Null

# ---------------------------------------------------------------------------
# 4) "a) " / "found" / " potential solution..." ->
#    single run "a) found potential solution to pair traveling and they meet goal of not being left with wrong pair "
# ---------------------------------------------------------------------------
Merge-WholeParagraphText "potential solution" "a) found potential solution to pair traveling and they meet goal of not being left with wrong pair " | Out-Null

# ---------------------------------------------------------------------------
# 5) "b) " / "tried" / " to take one at a time..." ->
#    single run "b) tried to take one at a time but it would leave impossible pairing on either side. "
# ---------------------------------------------------------------------------
Merge-WholeParagraphText "take one at a time" "b) tried to take one at a time but it would leave impossible pairing on either side. " | Out-Null

# ---------------------------------------------------------------------------
# 6) "b)" stays its own run; " " / "his" / " solution to find the number..." ->
#    single run " his solution to find the number in each sock will work for all cases and all colors."
# ---------------------------------------------------------------------------
$mergedHis = " his solution to find the number in each sock will work for all cases and all colors."
Merge-WholeParagraphText "solution to find the number" ("b)" + $mergedHis) | Out-Null
Split-RunAt "solution to find the number" 2 | Out-Null

# ---------------------------------------------------------------------------
# 7) "a)" and " The constraints " stay their own runs; "are you loose a finger..." /
#    "start" / " on one opposite finger..." ->
#    single run "are you loose a finger on a total count because u start on one opposite finger from both ways when counting "
# ---------------------------------------------------------------------------
$mergedStart = "are you loose a finger on a total count because u start on one opposite finger from both ways when counting "
Merge-WholeParagraphText "loose a finger" ("a) The constraints " + $mergedStart) | Out-Null
Split-RunAt "loose a finger" 19 | Out-Null

# ---------------------------------------------------------------------------
# 8) Insert new paragraphs (problem 3 potential solution) after the
#    "b) The sub goals are to find at what count do you count all fingers..." paragraph
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("The sub goals are to find at what count")) {
        $anchor = $p.Range
        $anchor.Collapse(0)  # collapse to end (wdCollapseEnd = 0)

        $anchor.InsertParagraphAfter()
        $anchor.Collapse(0)
        $anchor.MoveStart(1, 1)
        $anchor.ParagraphFormat.LeftIndent = 18

        $anchor.InsertParagraphAfter()
        $anchor.Collapse(0)
        $anchor.MoveStart(1, 1)
        $anchor.ParagraphFormat.LeftIndent = 18
        $anchor.InsertAfter("3.  ")

        $anchor.InsertParagraphAfter()
        $anchor.Collapse(0)
        $anchor.MoveStart(1, 1)
        $anchor.ParagraphFormat.LeftIndent = 18
        $anchor.InsertAfter("a) The solution to the sub problem is to count from one form the first finger till you count all your fingers once landing on the middle finger ")

        $anchor.InsertParagraphAfter()
        $anchor.Collapse(0)
        $anchor.MoveStart(1, 1)
        $anchor.ParagraphFormat.LeftIndent = 18

        break
    }
}
